$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.3301107639316767
$ws.Cells.Item(2, 2).Value = 0.415611196524354
$ws.Cells.Item(2, 3).Value = 0.2445170188243952
$ws.Cells.Item(2, 4).Value = 0.3903507817387732
$ws.Cells.Item(2, 5).Value = 0.2677516320288421
$ws.Cells.Item(2, 11).Value = 2.310775347521737
$ws.Cells.Item(2, 12).Value = 2.909278375670477
$ws.Cells.Item(2, 13).Value = 1.711619131770767
$ws.Cells.Item(2, 14).Value = 2.732455472171412
$ws.Cells.Item(2, 15).Value = 1.874261424201894
$ws.Cells.Item(2, 16).Value = 67.27676
$ws.Cells.Item(2, 17).Value = 106.5411467503276
$ws.Cells.Item(2, 18).Value = 44.49989721173976
$ws.Cells.Item(2, 19).Value = 88.79743446678775
$ws.Cells.Item(2, 20).Value = 48.77994374777913
$ws.Cells.Item(2, 21).Value = 0.2011838440952395
$ws.Cells.Item(2, 22).Value = 0.288744343116168
$ws.Cells.Item(2, 23).Value = 0.1016229933579305
$ws.Cells.Item(2, 24).Value = 0.2663293497599659
$ws.Cells.Item(2, 25).Value = 0.1313573526171028
$ws.Cells.Item(2, 26).Value = 0.8466649869005678
$ws.Cells.Item(2, 27).Value = 0.934499166025841
$ws.Cells.Item(2, 28).Value = 0.6967082179940818
$ws.Cells.Item(2, 29).Value = 0.9184631542161314
$ws.Cells.Item(2, 30).Value = 0.7580821531328032
$ws.Cells.Item(3, 6).Value = 7.003892718045113
$ws.Cells.Item(3, 7).Value = 9.852527728921453
$ws.Cells.Item(3, 8).Value = 4.14986917649107
$ws.Cells.Item(3, 9).Value = 9.015014774712206
$ws.Cells.Item(3, 10).Value = 4.933294032909545
$ws.Cells.Item(3, 11).Value = 2.311284596954888
$ws.Cells.Item(3, 12).Value = 3.25133415054408
$ws.Cells.Item(3, 13).Value = 1.369456828242053
$ws.Cells.Item(3, 14).Value = 2.974954875655028
$ws.Cells.Item(3, 15).Value = 1.627987030860151
$ws.Cells.Item(3, 16).Value = 68.48976
$ws.Cells.Item(3, 17).Value = 108.4255775396452
$ws.Cells.Item(3, 18).Value = 53.60549134344187
$ws.Cells.Item(3, 19).Value = 84.10651852181266
$ws.Cells.Item(3, 20).Value = 56.27728159932783
$ws.Cells.Item(3, 21).Value = 0.1968441701163827
$ws.Cells.Item(3, 22).Value = 0.3289385740419548
$ws.Cells.Item(3, 23).Value = 0.03927090336922618
$ws.Cells.Item(3, 24).Value = 0.297725408684395
$ws.Cells.Item(3, 25).Value = 0.0871737810077515
$ws.Cells.Item(3, 26).Value = 0.8160507378129342
$ws.Cells.Item(3, 27).Value = 0.9602856018097323
$ws.Cells.Item(3, 28).Value = 0.4935148305320233
$ws.Cells.Item(3, 29).Value = 0.934966440598985
$ws.Cells.Item(3, 30).Value = 0.6508775856717317
$ws.Cells.Item(4, 1).Value = 0.3297495105325244
$ws.Cells.Item(4, 2).Value = 0.4154772426198125
$ws.Cells.Item(4, 3).Value = 0.2446516004416382
$ws.Cells.Item(4, 4).Value = 0.3900524433988688
$ws.Cells.Item(4, 5).Value = 0.2675365285236319
$ws.Cells.Item(4, 6).Value = 7.002348170626321
$ws.Cells.Item(4, 7).Value = 9.852903518571656
$ws.Cells.Item(4, 8).Value = 4.14871823761041
$ws.Cells.Item(4, 9).Value = 9.018141007762488
$ws.Cells.Item(4, 10).Value = 4.929829151315143
$ws.Cells.Item(4, 11).Value = 2.308927180086172
$ws.Cells.Item(4, 12).Value = 3.730927800177077
$ws.Cells.Item(4, 13).Value = 1.194623138859719
$ws.Cells.Item(4, 14).Value = 3.033289829101107
$ws.Cells.Item(4, 15).Value = 1.577356635019926
$ws.Cells.Item(4, 16).Value = 70.30318
$ws.Cells.Item(4, 17).Value = 149.3804346691571
$ws.Cells.Item(4, 18).Value = 15.87156089226574
$ws.Cells.Item(4, 19).Value = 92.25890347765299
$ws.Cells.Item(4, 20).Value = 47.7981555424754
$ws.Cells.Item(4, 21).Value = 0.1938043747894106
$ws.Cells.Item(4, 22).Value = 0.3763101445064044
$ws.Cells.Item(4, 23).Value = 0.01207763732108842
$ws.Cells.Item(4, 24).Value = 0.306968598326241
$ws.Cells.Item(4, 25).Value = 0.07649327726257089
$ws.Cells.Item(4, 26).Value = 0.7891922249792941
$ws.Cells.Item(4, 27).Value = 1.005936367885767
$ws.Cells.Item(4, 28).Value = 0.2779109044429598
$ws.Cells.Item(4, 29).Value = 0.9210366299083371
$ws.Cells.Item(4, 30).Value = 0.5798096120396706
